# "Minor amends in pfscompd"
# - flips two existing CNSR flags (B40, B41) from 0 to 1
# - appends 18 new observations (rows 42-59) to the pfscompd sheet
# - makes pfscompd the active/selected sheet (was oscompd), with B58 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pfscompd")

# Amend two pre-existing rows: CNSR goes from 0 to 1
$ws.Range("B40").Value = 1
$ws.Range("B41").Value = 1

# New rows appended at the bottom of the table (AVAL, CNSR, TRT=Competimab)
$newRows = @(
    @(42, 5.3666666597222203, 1),
    @(43, 5.9986111041666597, 1),
    @(44, 6.6305555486111096, 1),
    @(45, 7.2624999930555498, 0),
    @(46, 7.8944444374999998, 0),
    @(47, 8.52638888194444,   0),
    @(48, 9.15833332638889,   0),
    @(49, 9.7902777708333293, 0),
    @(50, 10.422222215277801, 0),
    @(51, 11.054166659722201, 0),
    @(52, 11.686111104166701, 1),
    @(53, 12.318055548611101, 1),
    @(54, 12.9499999930556,   1),
    @(55, 13.581944437500001, 1),
    @(56, 14.213888881944399, 1),
    @(57, 14.845833326388901, 1),
    @(58, 15.477777770833301, 0),
    @(59, 16.109722215277699, 0)
)

foreach ($row in $newRows) {
    $rowNum = $row[0]
    $ws.Cells.Item($rowNum, 1).Value = $row[1]
    $ws.Cells.Item($rowNum, 2).Value = $row[2]
    $ws.Cells.Item($rowNum, 3).Value = "Competimab"
}

# Switch the active tab to pfscompd (oscompd was previously active/selected)
# and leave the selection on the last entered cell, like a user scrolling
# down the sheet after typing in the new rows.
$ws.Activate() | Out-Null
$ws.Range("B58").Select() | Out-Null
